$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Move the "Aufwand:" label block from L6:M6 to I16:J16
# ---------------------------------------------------------------------------
$ws.Range("L6:M6").Copy($ws.Range("I16"))
$ws.Range("L6:M6").Clear()

# ---------------------------------------------------------------------------
# 2. Move the "Text zum Diagramm:" block from L9:M9 to G20:H20
# ---------------------------------------------------------------------------
$ws.Range("L9:M9").Copy($ws.Range("G20"))
$ws.Range("L9:M9").Clear()

# ---------------------------------------------------------------------------
# 3. Extend the "T(N) = ..." rich text (now in J16) with a superscript "k",
#    a subscript "i=1" and a separating space right after the first
#    superscript "2" (i.e. after "(N/2)^2").
# ---------------------------------------------------------------------------
$tn = $ws.Range("J16")
$tn.Characters(14, 0).Text = "ki=1 "
$tn.Characters(13, 1).Font.Superscript = $true
$tn.Characters(14, 1).Font.Superscript = $true
$tn.Characters(15, 3).Font.Subscript = $true
$tn.Characters(22, 1).Font.Superscript = $true
$tn.Characters(29, 1).Font.Superscript = $true

# ---------------------------------------------------------------------------
# 4. Append a sentence about logarithmic scaling to the chart description
#    (now in H20).
# ---------------------------------------------------------------------------
$desc = $ws.Range("H20")
$desc.Characters(182, 0).Text = "Die Skalierung ist bei beiden AchsenLogarithmisch."
$desc.Characters(108, 13).Font.Name = "Calibri"

# ---------------------------------------------------------------------------
# 5. Row height adjustments
# ---------------------------------------------------------------------------
$ws.Rows(6).AutoFit()
$ws.Rows(9).RowHeight = 15.75
$ws.Rows(16).RowHeight = 17.25
$ws.Rows(20).RowHeight = 99.75

# ---------------------------------------------------------------------------
# 6. Column width adjustments
# ---------------------------------------------------------------------------
$ws.Columns(7).ColumnWidth = 22.022135416666668
$ws.Columns(8).ColumnWidth = 41.877604166666664
$ws.Columns(10).ColumnWidth = 23.736979166666668

# ---------------------------------------------------------------------------
# 7. Resize / reposition the chart and adjust its inner plot-area layout
# ---------------------------------------------------------------------------
$co = $ws.ChartObjects(1)
$co.Left = 167.87492125984252
$co.Top = 7.5
$co.Width = 424.1250787401575
$co.Height = 263.25

$chart = $co.Chart
$pa = $chart.PlotArea
$pa.Left = 0.12348373331050315
$pa.Top = 0.11443509732223643

# ---------------------------------------------------------------------------
# 8. Selection
# ---------------------------------------------------------------------------
$ws.Range("H20").Select()
